$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "What are the causes of World War II? Explain.",
    "Why did Nazi Germany and Imperial Japan surrender in 1945?",
    "When did the Bangladesh Liberation War start?",
    "Write   two   basic   difference   between   Gentlemen   Opium   Monopoly   and   BengalGovernment Opium Monopoly.",
    "Which one was the black gold of India?Why?",
    "Why the battle of palashi was called a fateful battle?",
    "When did Bangladesh gain independence?",
    "What did East Pakistan change its name to in 1971?",
    "What is the difference between an enzyme and substrate?",
    "Why is DNA replication necessary?",
    "Define geography and environmental geography",
    "Discuss the various types of region and their functions in Bangladesh",
    "Discuss the characteristics of agriculture in Bangladesh",
    "What do you understand by arsenic pollution in ground water? 2",
    "Discuss the main greenhouses gases (GHG) and their effects on environment?",
    "What are the characteristics and distinctions between gnosticism, epicureanism, and stoicism?",
    "Why does process philosophy need subjective aims?",
    "What did Epicurus say about virtue?",
    "Why did Greek philosophers only have first names?"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$styledRange = $ws.Range("A10:A20")
$styledRange.Font.Name = "Arial"
$styledRange.Font.Size = 10
$styledRange.Font.Color = 0

$ws.Range("A20").Select()
